# Add data for 2022-10-09
# - Rename sheet/tab title from "Through 2022-09-29" to "Through 2022-10-01"
# - Rename the "September (through 09-29)" month label to plain "September"
# - Insert a new "October (through 10-01)" row with partial data (D:I only)
# - Update the September row values and the Total row values accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-10-01"

# September row (row 10): drop the "(through 09-29)" qualifier, update counts.
$ws.Range("A10").Value = "September"
$ws.Range("B10").Value = 32
$ws.Range("C10").Value = 48
$ws.Range("D10").Value = 76
$ws.Range("E10").Value = 58
$ws.Range("F10").Value = 73
$ws.Range("G10").Value = 117
$ws.Range("H10").Value = 178
$ws.Range("I10").Value = 148

# Insert a new row for October before the existing Total row so Total stays last.
$ws.Rows("11").Insert()

# October row (row 11): only partial-month data (no 2015/2016 columns).
$ws.Range("A11").Value = "October (through 10-01)"
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 8
$ws.Range("I11").Value = 5

# Copy the style from the month label cell above onto the new October label.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) # xlPasteFormats

# Total row is now row 12; update with the new totals.
$ws.Range("B12").Value = 226
$ws.Range("C12").Value = 429
$ws.Range("D12").Value = 633
$ws.Range("E12").Value = 551
$ws.Range("F12").Value = 423
$ws.Range("G12").Value = 905
$ws.Range("H12").Value = 1255
$ws.Range("I12").Value = 1287

# Column A is slightly narrower now that "(through 09-29)" text is gone.
# (ColumnWidth is pixel-quantized by Excel; 22.8 is the nearest settable
# value that rounds to the target ~23.71 width.)
$ws.Columns("A").ColumnWidth = 22.8
